$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.361.79'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.936.03'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '0.7509'
$ws.Range("E5").Value = '  +5.42%  '
$ws.Range("D6").Value = '243.95'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '27.95'
$ws.Range("E8").Value = '  +2.34%  '
$ws.Range("D9").Value = '0.3186'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '0.07033'
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7830'
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").Value = '0.08043'
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").Value = '1.924.13'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").Value = '5.416'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").Value = '93.39'
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").Value = '14.46'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '30.352.41'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").Value = '6.069'
$ws.Range("E18").Value = '  +5.54%  '
$ws.Range("D19").Value = '252.69'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007992'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '2.184.20'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '6.712'
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").Value = '9.541'
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("D26").Value = '164.52'
$ws.Range("E26").Value = '  -0.92%  '
$ws.Range("D27").Value = '19.12'
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1310'
$ws.Range("E28").Value = '  +3.28%  '
$ws.Range("D29").Value = '2.229'
$ws.Range("E29").Value = '  -1.71%  '
$ws.Range("D30").Value = '1.377'
$ws.Range("E30").Value = '  +1.52%  '
$ws.Range("D31").Value = '1.532'
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").Value = '4.419'
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = '4.144'
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.343'
$ws.Range("E34").Value = '  +6.31%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.05247'
$ws.Range("E35").Value = '  +2.00%  '
$ws.Range("D36").Value = '0.7577'
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("D37").Value = '2.792'
$ws.Range("E37").Value = '  +0.89%  '
$ws.Range("D38").Value = '0.01959'
$ws.Range("E38").Value = '  +0.62%  '
$ws.Range("D39").Value = '2.814'
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.10'
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.565'
$ws.Range("E41").Value = '  +3.90%  '
$ws.Range("D42").Value = '0.4507'
$ws.Range("E42").Value = '  +0.81%  '
$ws.Range("D43").Value = '1.981'
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8380'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").Value = '10.04'
$ws.Range("E46").Value = '  +3.74%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.689'
$ws.Range("E47").Value = '  +3.56%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '101.63'
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("D49").Value = '37.85'
$ws.Range("E49").Value = '  +4.17%  '
$ws.Range("D50").Value = '980.85'
$ws.Range("E50").Value = '  +7.52%  '
$ws.Range("D51").Value = '0.1232'
$ws.Range("E51").Value = '  +9.26%  '
